# Apply WRI input-data update to the "About" sheet:
#  - insert a new row above the old row 10, pushing the "Methodology" section
#    (and everything below it) down by one row
#  - the previously-blank row 9 now holds a hyperlinked citation to a cached
#    web-archive copy of the IEA source document (the original link having
#    gone stale)
#  - the newly-inserted row 10 holds an explanatory note about the cached link

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

$archiveUrl = "https://web.archive.org/web/20170918224035/http://www.iea.org/publications/freepublications/publication/transport2009.pdf"

# Shift "Methodology" (old row 10) and everything below it down by one row.
$ws.Rows.Item(10).Insert()

# Row 9: hyperlink to the cached web-archive copy of the source PDF.
$ws.Range("B9").Value = $archiveUrl
$ws.Hyperlinks.Add($ws.Range("B9"), $archiveUrl)
$ws.Range("B9").Style = "Hyperlink"

# Row 10: note explaining why the cached link is used.
$ws.Range("B10").Value = "(Cached web archive link - original link unavailable)"

# Restore the on-screen selection state recorded in the saved workbook.
$ws5 = $wb.Worksheets.Item("PCiCDTdtTDM")
[void]$ws5.Activate()
[void]$ws5.Range("B2").Select()

[void]$ws.Activate()
[void]$ws.Range("B10").Select()

Write-Output "About sheet updated with cached web-archive citation"
